$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and swap two row pairs) per
# the latest scrape: Tue Apr 23 05:39:00 UTC 2024

$ws.Range('D2').Value = '66.503.53'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '3.184.41'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.14'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.93'
$ws.Range('E6').Value = '  +3.08%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.182.84'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.550'
$ws.Range('E9').Value = '  +2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.69'
$ws.Range('E11').Value = '  -7.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.513'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000268'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.90'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').Value = '3.706.68'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = '66.493.71'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.44'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = '3.182.31'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '514.03'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.54'
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.736'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.27'
$ws.Range('E23').Value = '  +3.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.98'
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.73'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.01'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.24'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.42'
$ws.Range('E29').Value = '  +8.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.12'
$ws.Range('E30').Value = '  +8.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.10'
$ws.Range('E31').Value = '  +4.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.16'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('B33').Value = 'Mantle'
$ws.Range('C33').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.56'
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '513.38'
$ws.Range('E36').Value = '  +5.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.73'
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0897'
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0424'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.05'
$ws.Range('E40').Value = '  +2.12%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.126'
$ws.Range('E41').Value = '  +5.87%  '
$ws.Range('D42').Value = '0.0₃0686'
$ws.Range('E42').Value = '  +6.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.302'
$ws.Range('E43').Value = '  +3.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.85'
$ws.Range('E44').Value = '  -4.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '2.852.66'
$ws.Range('E46').Value = '  -5.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.51'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  +4.14%  '
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.61'
$ws.Range('E51').Value = '  +7.53%  '
